$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.801.20"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.55%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.817.31"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.42%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "600.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.11%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.36"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.22%  "

$ws.Range("E7").Value = "  +0.09%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.519"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.20%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.160"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.75%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.29"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.82%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.451"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.82%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000253"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.49%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "35.84"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.89%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.458.45"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.42%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.800.45"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.11%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.808.93"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.53%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.41"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.04%  "

$ws.Range("E18").Value = "  +1.45%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "462.45"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.47%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.89"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.10%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.700"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.63%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.0000149"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.32%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.45"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.22%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.09"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.63%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.11"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.17%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.07"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.32%  "

$ws.Range("E28").Value = "  -0.02%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.966.53"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.39%  "

$ws.Range("E30").Value = "  +0.35%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.40"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.96%  "

$ws.Range("E32").Value = "  +2.32%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "29.60"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.21%  "

$ws.Range("E34").Value = "  -0.05%  "

$ws.Range("B35").Value = "Aptos"
$ws.Range("C35").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.09"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.13%  "

$ws.Range("B36").Value = "RenzoRestakedETH"
$ws.Range("C36").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.758.90"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.11%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.1000"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.02%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.29"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.11%  "

$ws.Range("E39").Value = "  -0.13%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.996"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.01%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.80"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.12%  "

$ws.Range("E42").Value = "  +0.01%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "48.11"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.38%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "28.75"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +11.82%  "

$ws.Range("B46").Value = "Arweave"
$ws.Range("C46").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "43.62"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.73%  "

$ws.Range("B47").Value = "TheGraph"
$ws.Range("C47").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.300"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.55%  "

$ws.Range("E48").Value = "  +11.57%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.35"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.37%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "148.65"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.16%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.84"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.12%  "
